$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.200.58'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.256.81'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +5.31%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '577.19'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.53'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +6.57%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.248.48'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +5.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +3.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.04'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +8.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.165'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +4.00%  '
$ws.Range("E12").Value = '  +3.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.97'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +3.28%  '
$ws.Range("E14").Value = '  +4.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.775.00'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +5.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '561.61'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +11.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.312.64'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +2.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.253.17'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +5.34%  '
$ws.Range("E19").Value = '  +2.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.17'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +5.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.54'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +4.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.748'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +6.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.86'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +7.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.66'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +5.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '82.63'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +3.15%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.39'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +17.84%  '
$ws.Range("E28").Value = '  +5.86%  '
$ws.Range("E29").Value = '  +5.52%  '
$ws.Range("E30").Value = '  +5.28%  '
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  +3.09%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '569.70'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +10.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +5.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0463'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +13.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.25'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0871'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +6.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.08'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +11.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.129'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.166.28'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +6.61%  '
$ws.Range("E43").Value = '  +1.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.277'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +9.94%  '
$ws.Range("E45").Value = '  +6.08%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.75'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +3.76%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₃0559'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '125.59'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("E50").Value = '  +2.51%  '
$ws.Range("E51").Value = '  +7.26%  '
